$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Column C values are written as text (matching original inlineStr type)
$ws.Range("C2:C47").NumberFormat = "@"

$ws.Range("B2").Value = "26-08-2025."
$ws.Range("C2").Value = "6"
$ws.Range("B3").Value = "02-08-2025."
$ws.Range("C3").Value = "3"
$ws.Range("B4").Value = "20-08-2025."
$ws.Range("C4").Value = "2"
$ws.Range("B5").Value = "12-08-2025."
$ws.Range("C5").Value = "7"
$ws.Range("B6").Value = "22-08-2025."
$ws.Range("C6").Value = "5"
$ws.Range("B7").Value = "05-08-2025."
$ws.Range("C7").Value = "4"
$ws.Range("B8").Value = "06-08-2025."
$ws.Range("C8").Value = "2"
$ws.Range("B9").Value = "18-08-2025."
$ws.Range("C9").Value = "6"
$ws.Range("B10").Value = "23-08-2025."
$ws.Range("C10").Value = "5"
$ws.Range("B11").Value = "09-08-2025."
$ws.Range("C11").Value = "5"
$ws.Range("B12").Value = "09-08-2025."
$ws.Range("C12").Value = "2"
$ws.Range("B13").Value = "07-08-2025."
$ws.Range("C13").Value = "2"
$ws.Range("B14").Value = "09-08-2025."
$ws.Range("C14").Value = "6"
$ws.Range("B15").Value = "11-08-2025."
$ws.Range("C15").Value = "6"
$ws.Range("B16").Value = "10-08-2025."
$ws.Range("C16").Value = "6"
$ws.Range("B17").Value = "12-08-2025."
$ws.Range("C17").Value = "5"
$ws.Range("B18").Value = "12-08-2025."
$ws.Range("C18").Value = "3"
$ws.Range("B19").Value = "12-08-2025."
$ws.Range("C19").Value = "2"
$ws.Range("B20").Value = "13-08-2025."
$ws.Range("C20").Value = "5"
$ws.Range("B21").Value = "09-08-2025."
$ws.Range("C21").Value = "11"
$ws.Range("B22").Value = "13-08-2025."
$ws.Range("C22").Value = "2"
$ws.Range("B23").Value = "28-08-2025."
$ws.Range("C23").Value = "6"
$ws.Range("B24").Value = "29-08-2025."
$ws.Range("C24").Value = "6"
$ws.Range("B25").Value = "15-08-2025."
$ws.Range("C25").Value = "2"
$ws.Range("B26").Value = "18-08-2025."
$ws.Range("C26").Value = "2"
$ws.Range("B27").Value = "16-08-2025."
$ws.Range("C27").Value = "5"
$ws.Range("B28").Value = "17-08-2025."
$ws.Range("C28").Value = "5"
$ws.Range("B29").Value = "04-08-2025."
$ws.Range("C29").Value = "2"
$ws.Range("B30").Value = "19-08-2025."
$ws.Range("C30").Value = "6"
$ws.Range("B31").Value = "24-08-2025."
$ws.Range("C31").Value = "5"
$ws.Range("B32").Value = "22-08-2025."
$ws.Range("C32").Value = "3"
$ws.Range("B33").Value = "22-08-2025."
$ws.Range("C33").Value = "2"
$ws.Range("B34").Value = "23-08-2025."
$ws.Range("C34").Value = "5"
$ws.Range("B35").Value = "23-08-2025."
$ws.Range("C35").Value = "3"
$ws.Range("B36").Value = "21-08-2025."
$ws.Range("C36").Value = "11"
$ws.Range("B37").Value = "30-08-2025."
$ws.Range("C37").Value = "6"
$ws.Range("B38").Value = "14-08-2025."
$ws.Range("C38").Value = "12"
$ws.Range("B39").Value = "11-08-2025."
$ws.Range("C39").Value = "12"
$ws.Range("B40").Value = "03-08-2025."
$ws.Range("C40").Value = "2"
$ws.Range("B41").Value = "08-08-2025."
$ws.Range("C41").Value = "6"
$ws.Range("B42").Value = "31-08-2025."
$ws.Range("C42").Value = "6"
$ws.Range("B43").Value = "27-08-2025."
$ws.Range("C43").Value = "6"
$ws.Range("B44").Value = "25-08-2025."
$ws.Range("C44").Value = "5"
$ws.Range("B45").Value = "05-08-2025."
$ws.Range("C45").Value = "11"
$ws.Range("B46").Value = "22-08-2025."
$ws.Range("C46").Value = "11"
$ws.Range("B47").Value = "21-08-2025."
$ws.Range("C47").Value = "6"
